$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number must be force-written as text
# (to match the source inlineStr cells) by assigning a text formula and then
# collapsing it to a static value via Copy + PasteSpecial(xlPasteValues = -4163).
# Plain .Value assignment would otherwise let Excel auto-convert the text into a
# number, which also pulls in a "Text" number-format style that is not in the diff.

$ws.Range("D2").Value = "57.582.23"
$ws.Range("E2").Value = "  -4.60%  "
$ws.Range("D3").Value = "2.936.05"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Formula = "=`"550.35`""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -4.02%  "
$ws.Range("D6").Formula = "=`"130.93`""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +3.90%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").Value = "2.932.27"
$ws.Range("D10").Formula = "=`"0.126`""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").Formula = "=`"4.76`""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -6.03%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Formula = "=`"0.0000221`""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Formula = "=`"32.92`""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "3.421.36"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Formula = "=`"6.89`""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +5.50%  "
$ws.Range("D18").Value = "2.934.04"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").Value = "57.659.97"
$ws.Range("E19").Value = "  -4.44%  "
$ws.Range("D20").Formula = "=`"417.78`""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Formula = "=`"6.99`""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").Formula = "=`"13.12`""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Formula = "=`"79.85`""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Formula = "=`"2.46`""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Formula = "=`"7.43`""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").Formula = "=`"25.21`""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").Formula = "=`"5.99`""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").Formula = "=`"5.68`""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").Formula = "=`"0.935`""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "0.0₃0692"
$ws.Range("E37").Value = "  +2.92%  "
$ws.Range("D38").Formula = "=`"48.25`""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("D39").Formula = "=`"8.76`""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("D40").Formula = "=`"2.55`""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Formula = "=`"376.37`""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Formula = "=`"0.0345`""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").Value = "2.700.59"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Formula = "=`"122.38`""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Formula = "=`"0.238`""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Formula = "=`"23.09`""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("E51").Value = "  -0.38%  "

$excel.CutCopyMode = 0

